$d = $word.ActiveDocument
$s = $d.Styles.Item("Normal")
$s.Font.Name = "Arial"
